# Append new data row 45 to the CityResaleNum sheet, per commit:
# "Realestate Update resale numbers 2024-01-10 21:25"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 45

# Columns A-D are text values. Setting Value directly would let Excel's
# auto-detection convert date/number-looking strings (e.g. "2024-01-10",
# "01") into a date serial / numeric value with an applied number format.
# Force them to Text format first, assign the string, then clear the
# formatting again so the cell ends up as a plain (unstyled) shared string,
# matching the rest of the sheet's data rows.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws.Cells.Item($row, 1) "2024-01-10"
Set-TextValue $ws.Cells.Item($row, 2) "21:25:15"
Set-TextValue $ws.Cells.Item($row, 3) "Wednesday"
Set-TextValue $ws.Cells.Item($row, 4) "01"

# Columns E-T are numeric values.
$ws.Cells.Item($row, 5).Value = 139498
$ws.Cells.Item($row, 6).Value = 142666
$ws.Cells.Item($row, 7).Value = 172169
$ws.Cells.Item($row, 8).Value = 147979
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 119292
$ws.Cells.Item($row, 11).Value = 224805
$ws.Cells.Item($row, 12).Value = 251492
$ws.Cells.Item($row, 13).Value = 185247
$ws.Cells.Item($row, 14).Value = 110449
$ws.Cells.Item($row, 15).Value = 40798
$ws.Cells.Item($row, 16).Value = 30883
$ws.Cells.Item($row, 17).Value = 72799
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42380
$ws.Cells.Item($row, 20).Value = -1
